$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H6").Value = 15023.75
$ws_ALC.Range("I6").Value = 2078.2
$ws_ALC.Range("J6").Value = 36599.668
$ws_ALC.Range("K6").Value = 6234.599999999999
$ws_ALC.Range("L6").Value = 109799.004
$ws_ALC.Range("M6").Value = -6122.599999999999
$ws_ALC.Range("N6").Value = -110023.004

$ws_ALC.Range("H17").Value = 750.14703
$ws_ALC.Range("J17").Value = 750.14703
$ws_ALC.Range("L17").Value = 2250.44109
$ws_ALC.Range("N17").Value = -2586.44109

$ws_ALC.Range("H108").Value = 37999.668
$ws_ALC.Range("J108").Value = 37999.668
$ws_ALC.Range("L108").Value = 37999.668
$ws_ALC.Range("N108").Value = -45679.668

$ws_ALC.Range("H132").Value = 24188.072
$ws_ALC.Range("I132").Value = 3511.9355
$ws_ALC.Range("J132").Value = 82457.17999999999
$ws_ALC.Range("K132").Value = 10535.8065
$ws_ALC.Range("L132").Value = 247371.54
$ws_ALC.Range("M132").Value = -8005.806500000001
$ws_ALC.Range("N132").Value = -252431.54

$ws_ALC.Range("H137").Value = 4884.6294
$ws_ALC.Range("I137").Value = 915.7222
$ws_ALC.Range("J137").Value = 12822.444
$ws_ALC.Range("K137").Value = 2747.1666
$ws_ALC.Range("L137").Value = 38467.33199999999
$ws_ALC.Range("M137").Value = -197.1666
$ws_ALC.Range("N137").Value = -43567.33199999999

$ws_ARM.Range("H61").Value = 2471
$ws_ARM.Range("I61").Value = 1884.4615
$ws_ARM.Range("K61").Value = 1884.4615
$ws_ARM.Range("M61").Value = -1672.4615

$ws_ARM.Range("H74").Value = 2247.5
$ws_ARM.Range("I74").Value = 1907.3
$ws_ARM.Range("J74").Value = 4799
$ws_ARM.Range("K74").Value = 1907.3
$ws_ARM.Range("L74").Value = 4799
$ws_ARM.Range("M74").Value = -1033.3
$ws_ARM.Range("N74").Value = -6547

$ws_ARM.Range("H77").Value = 2247.5
$ws_ARM.Range("I77").Value = 1907.3
$ws_ARM.Range("J77").Value = 4799
$ws_ARM.Range("K77").Value = 9536.5
$ws_ARM.Range("L77").Value = 23995
$ws_ARM.Range("M77").Value = -5168.5
$ws_ARM.Range("N77").Value = -32731

$ws_ARM.Range("H132").Value = 20001834
$ws_ARM.Range("I132").Value = 31251116
$ws_ARM.Range("J132").Value = 3110.889
$ws_ARM.Range("K132").Value = 93753348
$ws_ARM.Range("L132").Value = 9332.667000000001
$ws_ARM.Range("M132").Value = -93750818
$ws_ARM.Range("N132").Value = -14392.667

$ws_ARM.Range("H135").Value = 35809.668
$ws_ARM.Range("J135").Value = 35809.668
$ws_ARM.Range("L135").Value = 35809.668
$ws_ARM.Range("N135").Value = -45949.668

$ws_ARM.Range("H136").Value = 2471
$ws_ARM.Range("I136").Value = 1884.4615
$ws_ARM.Range("K136").Value = 5653.3845
$ws_ARM.Range("M136").Value = -3103.3845

$ws_BSM.Range("H99").Value = 2022.4688
$ws_BSM.Range("I99").Value = 1579.9584
$ws_BSM.Range("J99").Value = 3350
$ws_BSM.Range("K99").Value = 1579.9584
$ws_BSM.Range("L99").Value = 3350
$ws_BSM.Range("M99").Value = -81.95839999999998
$ws_BSM.Range("N99").Value = -6346

$ws_BSM.Range("H134").Value = 2686.5679
$ws_BSM.Range("I134").Value = 2037.674
$ws_BSM.Range("J134").Value = 3539.4
$ws_BSM.Range("K134").Value = 6113.022
$ws_BSM.Range("L134").Value = 10618.2
$ws_BSM.Range("M134").Value = -3578.022
$ws_BSM.Range("N134").Value = -15688.2

$ws_CRP.Range("H31").Value = 1436.99
$ws_CRP.Range("I31").Value = 716.2727
$ws_CRP.Range("J31").Value = 2836.0293
$ws_CRP.Range("K31").Value = 716.2727
$ws_CRP.Range("L31").Value = 2836.0293
$ws_CRP.Range("M31").Value = -421.2727
$ws_CRP.Range("N31").Value = -3426.0293

$ws_CRP.Range("H34").Value = 1436.99
$ws_CRP.Range("I34").Value = 716.2727
$ws_CRP.Range("J34").Value = 2836.0293
$ws_CRP.Range("K34").Value = 716.2727
$ws_CRP.Range("L34").Value = 2836.0293
$ws_CRP.Range("M34").Value = -514.2727
$ws_CRP.Range("N34").Value = -3240.0293

$ws_CRP.Range("H58").Value = 1462.6428
$ws_CRP.Range("I58").Value = 1191.6774
$ws_CRP.Range("J58").Value = 2226.2727
$ws_CRP.Range("K58").Value = 1191.6774
$ws_CRP.Range("L58").Value = 2226.2727
$ws_CRP.Range("M58").Value = -988.6774
$ws_CRP.Range("N58").Value = -2632.2727

$ws_CRP.Range("H60").Value = 10569.381
$ws_CRP.Range("I60").Value = 5500
$ws_CRP.Range("K60").Value = 5500
$ws_CRP.Range("M60").Value = -4989

$ws_CRP.Range("H132").Value = 60784.457
$ws_CRP.Range("I132").Value = 1060
$ws_CRP.Range("J132").Value = 144398.7
$ws_CRP.Range("K132").Value = 3180
$ws_CRP.Range("L132").Value = 433196.1
$ws_CRP.Range("M132").Value = -650
$ws_CRP.Range("N132").Value = -438256.1

$ws_CRP.Range("H134").Value = 390334.38
$ws_CRP.Range("I134").Value = 1044.0769
$ws_CRP.Range("J134").Value = 1402489.1
$ws_CRP.Range("K134").Value = 3132.2307
$ws_CRP.Range("L134").Value = 4207467.300000001
$ws_CRP.Range("M134").Value = -597.2307000000001
$ws_CRP.Range("N134").Value = -4212537.300000001

$ws_CRP.Range("H136").Value = 1462.6428
$ws_CRP.Range("I136").Value = 1191.6774
$ws_CRP.Range("J136").Value = 2226.2727
$ws_CRP.Range("K136").Value = 3575.0322
$ws_CRP.Range("L136").Value = 6678.8181
$ws_CRP.Range("M136").Value = -1025.0322
$ws_CRP.Range("N136").Value = -11778.8181

$ws_CUL.Range("H107").Value = 12869.125
$ws_CUL.Range("I107").Value = 11516.556
$ws_CUL.Range("J107").Value = 14608.143
$ws_CUL.Range("K107").Value = 34549.66800000001
$ws_CUL.Range("L107").Value = 43824.429
$ws_CUL.Range("M107").Value = -32629.66800000001
$ws_CUL.Range("N107").Value = -47664.429

$ws_CUL.Range("H119").Value = 3096.923
$ws_CUL.Range("I119").Value = 2098.375
$ws_CUL.Range("J119").Value = 4694.6
$ws_CUL.Range("K119").Value = 6295.125
$ws_CUL.Range("L119").Value = 14083.8
$ws_CUL.Range("M119").Value = -1457.125
$ws_CUL.Range("N119").Value = -23759.8

$ws_GSM.Range("H122").Value = 1943.6875
$ws_GSM.Range("I122").Value = 1918.1818
$ws_GSM.Range("K122").Value = 5754.5454
$ws_GSM.Range("M122").Value = -3304.5454

$ws_GSM.Range("H132").Value = 2830
$ws_GSM.Range("I132").Value = 2053.3635
$ws_GSM.Range("K132").Value = 6160.0905
$ws_GSM.Range("M132").Value = -3630.0905

$ws_LTW.Range("H40").Value = 4899.2104
$ws_LTW.Range("I40").Value = 3669.6155
$ws_LTW.Range("K40").Value = 3669.6155
$ws_LTW.Range("M40").Value = -3533.6155

$ws_LTW.Range("H122").Value = 35376.867
$ws_LTW.Range("I122").Value = 37796.645
$ws_LTW.Range("K122").Value = 113389.935
$ws_LTW.Range("M122").Value = -110939.935

$ws_LTW.Range("H132").Value = 3058.0625
$ws_LTW.Range("I132").Value = 2157.0588
$ws_LTW.Range("J132").Value = 4079.2
$ws_LTW.Range("K132").Value = 6471.176399999999
$ws_LTW.Range("L132").Value = 12237.6
$ws_LTW.Range("M132").Value = -3941.176399999999
$ws_LTW.Range("N132").Value = -17297.6

$ws_LTW.Range("H136").Value = 2123.48
$ws_LTW.Range("I136").Value = 1579.55
$ws_LTW.Range("K136").Value = 4738.65
$ws_LTW.Range("M136").Value = -2188.65

$ws_WVR.Range("H47").Value = 22034.5
$ws_WVR.Range("J47").Value = 22034.5
$ws_WVR.Range("L47").Value = 22034.5
$ws_WVR.Range("N47").Value = -23178.5

$ws_WVR.Range("H132").Value = 1418.56
$ws_WVR.Range("I132").Value = 1076.6316
$ws_WVR.Range("J132").Value = 2501.3333
$ws_WVR.Range("K132").Value = 3229.8948
$ws_WVR.Range("L132").Value = 7503.999899999999
$ws_WVR.Range("M132").Value = -699.8948
$ws_WVR.Range("N132").Value = -12563.9999

$ws_WVR.Range("H136").Value = 527732.2
$ws_WVR.Range("I136").Value = 667834.3
$ws_WVR.Range("J136").Value = 2349.25
$ws_WVR.Range("K136").Value = 2003502.9
$ws_WVR.Range("L136").Value = 7047.75
$ws_WVR.Range("M136").Value = -2000952.9
$ws_WVR.Range("N136").Value = -12147.75
